$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before DF. This shifts the old "nom" (DF) and
# "url_produit" (DG) columns one to the right (to DG and DH respectively),
# and leaves a fresh, empty DF column in their place.
$ws.Columns("DF").Insert()

# DF1 becomes the newest timestamp header cell (previously DE1 held
# "2026-02-01 16:15:30"; this commit adds the next scrape timestamp).
$ws.Range("DF1").Value = "2026-02-01 17:15:48"

# For every data row, the newly inserted price-history column simply
# repeats the most recent known price (the same value already stored in
# column DE, the previous last-scraped price column). Rows with no price
# recorded yet (DE holds an empty string) stay empty too.
$src = $ws.Range("DE2:DE206")
$dst = $ws.Range("DF2:DF206")
$dst.Value = $src.Value()
